$d = $word.ActiveDocument

# Locate the "{m:userdoc 'zone1'}" field text (currently stored as two runs:
# "{m" and ":userdoc 'zone1'}") and rewrite it as four separate runs:
# "{", "m", ":userdoc 'zone1'", "}" - matching the TokenIteratorFieldRewriterSplit
# output, without altering any other content in the document.

$full = $d.Content.Text
$needle = "{m:userdoc 'zone1'}"
$base = $full.IndexOf($needle)

if ($base -ge 0) {
    $pieces = @("{", "m", ":userdoc 'zone1'", "}")

    # Pass 1: carve out each piece's range using a same-length placeholder.
    # Changing the text (even to a throwaway value) is what makes the engine
    # split the run at this exact boundary and drop the stale w:rsidR carried
    # over from the original (now-merged) run.
    $offset = $base
    foreach ($piece in $pieces) {
        $len = $piece.Length
        $placeholder = "#".PadRight($len, "#")
        $r = $d.Range($offset, $offset + $len)
        $ft = $r.FormattedText
        $ft.Text = $placeholder
        $r.FormattedText = $ft
        $offset = $offset + $len
    }

    # Pass 2: write back the real text for each of the four ranges (same
    # positions/lengths as above, so offsets are still valid).
    $offset = $base
    foreach ($piece in $pieces) {
        $len = $piece.Length
        $r = $d.Range($offset, $offset + $len)
        $ft = $r.FormattedText
        $ft.Text = $piece
        $r.FormattedText = $ft
        $offset = $offset + $len
    }
}
